$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking values in C:G are stored as text (matches source data which
# uses inline strings for all line-item figures, including formatted numbers like "1,735").
$ws.Range("C2:G23").NumberFormat = "@"

$ws.Range("A2").Value = "  Cash & Equivalents"
$ws.Range("B2").Value = "cash_and_equivalents"
$ws.Range("C2").Value = "50"
$ws.Range("D2").Value = "50"
$ws.Range("E2").Value = "55"
$ws.Range("F2").Value = "61"
$ws.Range("G2").Value = "67"

$ws.Range("A3").Value = "  Accounts Receivable"
$ws.Range("B3").Value = "accounts_receivable"
$ws.Range("C3").Value = "100"
$ws.Range("D3").Value = "100"
$ws.Range("E3").Value = "112"
$ws.Range("F3").Value = "125"
$ws.Range("G3").Value = "140"

$ws.Range("A4").Value = "  Total Current Assets"
$ws.Range("B4").Value = "total_current_assets"
$ws.Range("C4").Value = "150"
$ws.Range("D4").Value = "150"
$ws.Range("E4").Value = "167"
$ws.Range("F4").Value = "186"
$ws.Range("G4").Value = "207"

$ws.Range("A5").Value = "  Property, Plant & Equipment"
$ws.Range("B5").Value = "property_plant_equipment"
$ws.Range("C5").Value = "300"
$ws.Range("D5").Value = "300"
$ws.Range("E5").Value = "315"
$ws.Range("F5").Value = "331"
$ws.Range("G5").Value = "347"

$ws.Range("A6").Value = "  Total Assets"
$ws.Range("B6").Value = "total_assets"
$ws.Range("C6").Value = "450"
$ws.Range("D6").Value = "450"
$ws.Range("E6").Value = "482"
$ws.Range("F6").Value = "517"
$ws.Range("G6").Value = "554"

$ws.Range("A7").Value = "  Accounts Payable"
$ws.Range("B7").Value = "accounts_payable"
$ws.Range("C7").Value = "80"
$ws.Range("D7").Value = "80"
$ws.Range("E7").Value = "86"
$ws.Range("F7").Value = "93"
$ws.Range("G7").Value = "101"

$ws.Range("A8").Value = "  Total Debt"
$ws.Range("B8").Value = "total_debt"
$ws.Range("C8").Value = "150"
$ws.Range("D8").Value = "150"
$ws.Range("E8").Value = "153"
$ws.Range("F8").Value = "156"
$ws.Range("G8").Value = "159"

$ws.Range("A9").Value = "  Total Liabilities"
$ws.Range("B9").Value = "total_liabilities"
$ws.Range("C9").Value = "230"
$ws.Range("D9").Value = "230"
$ws.Range("E9").Value = "239"
$ws.Range("F9").Value = "249"
$ws.Range("G9").Value = "260"

$ws.Range("A10").Value = "  Common Stock"
$ws.Range("B10").Value = "common_stock"
$ws.Range("C10").Value = "100"
$ws.Range("D10").Value = "100"
$ws.Range("E10").Value = "100"
$ws.Range("F10").Value = "100"
$ws.Range("G10").Value = "100"

$ws.Range("A11").Value = "  Prior Retained Earnings"
$ws.Range("B11").Value = "prior_retained_earnings"
$ws.Range("C11").Value = "100"
$ws.Range("D11").Value = "100"
$ws.Range("E11").Value = "100"
$ws.Range("F11").Value = "100"
$ws.Range("G11").Value = "100"

$ws.Range("A12").Value = "  Dividends"
$ws.Range("B12").Value = "dividends"
$ws.Range("C12").Value = "-10"
$ws.Range("D12").Value = "-10"
$ws.Range("E12").Value = "-10"
$ws.Range("F12").Value = "-11"
$ws.Range("G12").Value = "-12"

$ws.Range("A13").Value = "  Total Equity"
$ws.Range("B13").Value = "total_equity"
$ws.Range("C13").Value = "1,735"
$ws.Range("D13").Value = "3,005"
$ws.Range("E13").Value = "4,416"
$ws.Range("F13").Value = "6,084"
$ws.Range("G13").Value = "8,082"

$ws.Range("A14").Value = "  Total Liabilities & Equity"
$ws.Range("B14").Value = "total_liabilities_equity"
$ws.Range("C14").Value = "1,965"
$ws.Range("D14").Value = "3,235"
$ws.Range("E14").Value = "4,656"
$ws.Range("F14").Value = "6,333"
$ws.Range("G14").Value = "8,342"

$ws.Range("A15").Value = "  Revenue"
$ws.Range("B15").Value = "revenue"
$ws.Range("C15").Value = "1,000"
$ws.Range("D15").Value = "1,200"
$ws.Range("E15").Value = "1,320"
$ws.Range("F15").Value = "1,452"
$ws.Range("G15").Value = "1,597"

$ws.Range("A16").Value = "  Cost of Goods Sold"
$ws.Range("B16").Value = "cost_of_goods_sold"
$ws.Range("C16").Value = "-400"
$ws.Range("D16").Value = "-500"
$ws.Range("E16").Value = "-625"
$ws.Range("F16").Value = "-781"
$ws.Range("G16").Value = "-977"

$ws.Range("A17").Value = "  Gross Profit"
$ws.Range("B17").Value = "gross_profit"
$ws.Range("C17").Value = "1,400"
$ws.Range("D17").Value = "1,700"
$ws.Range("E17").Value = "1,945"
$ws.Range("F17").Value = "2,233"
$ws.Range("G17").Value = "2,574"

$ws.Range("A18").Value = "  Operating Expenses"
$ws.Range("B18").Value = "operating_expenses"
$ws.Range("C18").Value = "-100"
$ws.Range("D18").Value = "-120"
$ws.Range("E18").Value = "-130"
$ws.Range("F18").Value = "-137"
$ws.Range("G18").Value = "-158"

$ws.Range("A19").Value = "  Operating Income"
$ws.Range("B19").Value = "operating_income"
$ws.Range("C19").Value = "1,500"
$ws.Range("D19").Value = "1,820"
$ws.Range("E19").Value = "2,075"
$ws.Range("F19").Value = "2,371"
$ws.Range("G19").Value = "2,732"

$ws.Range("A20").Value = "  Interest Expense"
$ws.Range("B20").Value = "interest_expense"
$ws.Range("C20").Value = "-50"
$ws.Range("D20").Value = "-60"
$ws.Range("E20").Value = "-63"
$ws.Range("F20").Value = "-66"
$ws.Range("G20").Value = "-69"

$ws.Range("A21").Value = "  Income Before Tax"
$ws.Range("B21").Value = "income_before_tax"
$ws.Range("C21").Value = "1,550"
$ws.Range("D21").Value = "1,880"
$ws.Range("E21").Value = "2,138"
$ws.Range("F21").Value = "2,437"
$ws.Range("G21").Value = "2,801"

$ws.Range("A22").Value = "  Income Tax"
$ws.Range("B22").Value = "income_tax"
$ws.Range("C22").Value = "-75"
$ws.Range("D22").Value = "-90"
$ws.Range("E22").Value = "-108"
$ws.Range("F22").Value = "-130"
$ws.Range("G22").Value = "-156"

$ws.Range("A23").Value = "  Net Income"
$ws.Range("B23").Value = "net_income"
$ws.Range("C23").Value = "1,525"
$ws.Range("D23").Value = "1,850"
$ws.Range("E23").Value = "2,116"
$ws.Range("F23").Value = "2,429"
$ws.Range("G23").Value = "2,799"
